# "conditional parallel run added"
# Adds a new column D ("parallel_run" / "yes") to the "configuration" sheet,
# mirroring the look of the existing "headless" column (C).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("configuration")

# Copy the header/data formatting from column C (headless) into column D
# so the new column matches the existing table style (bold header, borders,
# wrap text, ...).
$ws.Range("C1").Copy() | Out-Null
$ws.Range("D1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

$ws.Range("C2").Copy() | Out-Null
$ws.Range("D2").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

# New header + value
$ws.Cells.Item(1, 4).Value = "parallel_run"
$ws.Cells.Item(2, 4).Value = "yes"

# Match the new column's width / header row height to the source workbook.
$ws.Columns.Item(4).ColumnWidth = 16.33
$ws.Rows.Item(1).RowHeight = 31.5

# Selection ends up on the newly added cell.
$ws.Range("D2").Select() | Out-Null
